$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 37 (duplicate "Flood Monitoring / National Framework for
# NMS/NHS Services" row whose question text duplicated a meteorological-law
# question); remaining rows shift up.
$ws.Rows(37).Delete()

# Restore the selection/scroll state recorded for the sheet after this edit
# (user had scrolled down and selected the (now last) data row, row 37).
$ws.Range("A37:XFD37").Select()
